$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.942.73"
$ws.Range("E2").Value = "  +3.55%  "
$ws.Range("D3").Value = "3.030.15"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.67%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.025.16"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +17.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.41%  "
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("E13").Value = "  +4.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.89%  "
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "3.538.88"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.05%  "
$ws.Range("D18").Value = "62.902.44"
$ws.Range("E18").Value = "  +3.30%  "
$ws.Range("D19").Value = "3.030.66"
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "452.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.698"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.82%  "
$ws.Range("E27").Value = "  +5.03%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.21%  "
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.62"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("E34").Value = "  +4.15%  "
$ws.Range("D35").Value = "0.0₃0862"
$ws.Range("E35").Value = "  +8.28%  "
$ws.Range("E36").Value = "  +3.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.03%  "
$ws.Range("E39").Value = "  +9.54%  "
$ws.Range("E40").Value = "  +4.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("E43").Value = "  +17.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +16.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "392.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0360"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.42%  "
$ws.Range("D47").Value = "2.721.06"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "132.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.94%  "
$ws.Range("E51").Value = "  +8.73%  "
